# Mise à jour du suivi de projet
#
# - Insère une nouvelle ligne d'activité (lecture doc adaptateur CSE-H53N) au
#   dessus de la liste des dates, ce qui décale toutes les dates existantes
#   (lignes 5 à 29) d'une ligne vers le bas (lignes 6 à 30) et ajoute une
#   nouvelle date (45828) en bas de tableau.
# - Les hauteurs de ligne explicites sur les lignes 2 à 4 sont supprimées.
# - La fusion de cellules A2:A4 est supprimée.
# - Les hyperliens sur C3/C4 sont supprimés (le texte affiché est conservé).
# - La sélection active passe à E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Décaler les dates existantes (lignes 5-29) d'une ligne vers le bas -----
# On copie d'abord le format (style de date) de la dernière ligne existante
# vers la nouvelle dernière ligne, puis on recopie les valeurs de bas en haut
# pour ne pas écraser de données avant de les avoir lues.
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
for ($r = 29; $r -ge 5; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(5, 1).Value = ""

# --- Nouvelle ligne 5 : activité "adaptateur RS232/RJ45 CSE-H53N" ----------
$ws.Range("B5").Value = "Lecture de la documentation de l'adaptateur RS232 à RJ45 ""CSE-H53N"" par Sollae"
$ws.Range("C5").Value = "CSE-H53N | RS-232 Serial To Ethernet Converter"

# --- Nettoyage de la mise en forme des lignes 2 à 4 -------------------------
# Hauteur de ligne automatique (supprime les ht="30"/"45" explicites)
$ws.Rows("2:4").AutoFit()

# Suppression de la fusion A2:A4
$ws.Range("A2:A4").UnMerge()

# Suppression des hyperliens (le texte visible C3/C4 reste inchangé)
$ws.Range("C3:C4").Hyperlinks.Delete()

# --- Sélection active --------------------------------------------------------
$ws.Range("E9").Select()
